$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 56.05275600000001
$ws.Range("H2").Value = 168.158268
$ws.Range("I2").Value = 0.06617112780234141
$ws.Range("J2").Value = 0.06617112780234141
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.166450999999999
$ws.Range("N2").Value = 12.499353
$ws.Range("O2").Value = 0.7126954333415383
$ws.Range("P2").Value = 0.7126954333415383
$ws.Range("Q2").Value = 233.541061288956
$ws.Range("R2").Value = 2101.869551600604
$ws.Range("S2").Value = 0.04715986060378802
$ws.Range("T2").Value = 0.04715986060378802

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 56.05275600000001
$ws.Range("H3").Value = 168.158268
$ws.Range("I3").Value = 0.06617112780234141
$ws.Range("J3").Value = 0.06617112780234141
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6655859999999999
$ws.Range("N3").Value = 1.996758
$ws.Range("O3").Value = 0.1138523176430159
$ws.Range("P3").Value = 0.1138523176430159
$ws.Range("Q3").Value = 37.307929655016
$ws.Range("R3").Value = 335.771366895144
$ws.Range("S3").Value = 0.007533736261348773
$ws.Range("T3").Value = 0.007533736261348773

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 56.05275600000001
$ws.Range("H4").Value = 168.158268
$ws.Range("I4").Value = 0.06617112780234141
$ws.Range("J4").Value = 0.06617112780234141
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.01401
$ws.Range("N4").Value = 3.04203
$ws.Range("O4").Value = 0.1734522490154458
$ws.Range("P4").Value = 0.1734522490154458
$ws.Range("Q4").Value = 56.83805511156002
$ws.Range("R4").Value = 511.5424960040401
$ws.Range("S4").Value = 0.01147753093720461
$ws.Range("T4").Value = 0.01147753093720461

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 238.00471
$ws.Range("H5").Value = 714.01413
$ws.Range("I5").Value = 0.2809681665424124
$ws.Range("J5").Value = 0.2809681665424124
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.166450999999999
$ws.Range("N5").Value = 12.499353
$ws.Range("O5").Value = 0.7126954333415383
$ws.Range("P5").Value = 0.7126954333415383
$ws.Range("Q5").Value = 991.6349619842099
$ws.Range("R5").Value = 8924.71465785789
$ws.Range("S5").Value = 0.2002447292091221
$ws.Range("T5").Value = 0.2002447292091221

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 238.00471
$ws.Range("H6").Value = 714.01413
$ws.Range("I6").Value = 0.2809681665424124
$ws.Range("J6").Value = 0.2809681665424124
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6655859999999999
$ws.Range("N6").Value = 1.996758
$ws.Range("O6").Value = 0.1138523176430159
$ws.Range("P6").Value = 0.1138523176430159
$ws.Range("Q6").Value = 158.41260291006
$ws.Range("R6").Value = 1425.71342619054
$ws.Range("S6").Value = 0.03198887694476251
$ws.Range("T6").Value = 0.03198887694476251

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 238.00471
$ws.Range("H7").Value = 714.01413
$ws.Range("I7").Value = 0.2809681665424124
$ws.Range("J7").Value = 0.2809681665424124
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.01401
$ws.Range("N7").Value = 3.04203
$ws.Range("O7").Value = 0.1734522490154458
$ws.Range("P7").Value = 0.1734522490154458
$ws.Range("Q7").Value = 241.3391559871
$ws.Range("R7").Value = 2172.0524038839
$ws.Range("S7").Value = 0.04873456038852777
$ws.Range("T7").Value = 0.04873456038852776

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 126.0396663333333
$ws.Range("H8").Value = 378.118999
$ws.Range("I8").Value = 0.1487917359336884
$ws.Range("J8").Value = 0.1487917359336884
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.166450999999999
$ws.Range("N8").Value = 12.499353
$ws.Range("O8").Value = 0.7126954333415383
$ws.Range("P8").Value = 0.7126954333415383
$ws.Range("Q8").Value = 525.138093834183
$ws.Range("R8").Value = 4726.242844507648
$ws.Range("S8").Value = 0.1060431907188998
$ws.Range("T8").Value = 0.1060431907188998

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 126.0396663333333
$ws.Range("H9").Value = 378.118999
$ws.Range("I9").Value = 0.1487917359336884
$ws.Range("J9").Value = 0.1487917359336884
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.6655859999999999
$ws.Range("N9").Value = 1.996758
$ws.Range("O9").Value = 0.1138523176430159
$ws.Range("P9").Value = 0.1138523176430159
$ws.Range("Q9").Value = 83.89023735613799
$ws.Range("R9").Value = 755.012136205242
$ws.Range("S9").Value = 0.01694028398217803
$ws.Range("T9").Value = 0.01694028398217803

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 126.0396663333333
$ws.Range("H10").Value = 378.118999
$ws.Range("I10").Value = 0.1487917359336884
$ws.Range("J10").Value = 0.1487917359336884
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.01401
$ws.Range("N10").Value = 3.04203
$ws.Range("O10").Value = 0.1734522490154458
$ws.Range("P10").Value = 0.1734522490154458
$ws.Range("Q10").Value = 127.8054820586634
$ws.Range("R10").Value = 1150.24933852797
$ws.Range("S10").Value = 0.02580826123261058
$ws.Range("T10").Value = 0.02580826123261058

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 145.0922953333333
$ws.Range("H11").Value = 435.276886
$ws.Range("I11").Value = 0.1712836531648339
$ws.Range("J11").Value = 0.1712836531648339
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.166450999999999
$ws.Range("N11").Value = 12.499353
$ws.Range("O11").Value = 0.7126954333415383
$ws.Range("P11").Value = 0.7126954333415383
$ws.Range("Q11").Value = 604.5199389838619
$ws.Range("R11").Value = 5440.679450854757
$ws.Range("S11").Value = 0.122073077416633
$ws.Range("T11").Value = 0.122073077416633

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 145.0922953333333
$ws.Range("H12").Value = 435.276886
$ws.Range("I12").Value = 0.1712836531648339
$ws.Range("J12").Value = 0.1712836531648339
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.6655859999999999
$ws.Range("N12").Value = 1.996758
$ws.Range("O12").Value = 0.1138523176430159
$ws.Range("P12").Value = 0.1138523176430159
$ws.Range("Q12").Value = 96.57140048173198
$ws.Range("R12").Value = 869.1426043355879
$ws.Range("S12").Value = 0.01950104088717883
$ws.Range("T12").Value = 0.01950104088717883

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 145.0922953333333
$ws.Range("H13").Value = 435.276886
$ws.Range("I13").Value = 0.1712836531648339
$ws.Range("J13").Value = 0.1712836531648339
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.01401
$ws.Range("N13").Value = 3.04203
$ws.Range("O13").Value = 0.1734522490154458
$ws.Range("P13").Value = 0.1734522490154458
$ws.Range("Q13").Value = 147.1250383909534
$ws.Range("R13").Value = 1324.12534551858
$ws.Range("S13").Value = 0.02970953486102203
$ws.Range("T13").Value = 0.02970953486102202

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 151.033905
$ws.Range("H14").Value = 453.101715
$ws.Range("I14").Value = 0.1782978134070997
$ws.Range("J14").Value = 0.1782978134070997
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 4.166450999999999
$ws.Range("N14").Value = 12.499353
$ws.Range("O14").Value = 0.7126954333415383
$ws.Range("P14").Value = 0.7126954333415383
$ws.Range("Q14").Value = 629.2753645211549
$ws.Range("R14").Value = 5663.478280690394
$ws.Range("S14").Value = 0.1270720373900217
$ws.Range("T14").Value = 0.1270720373900217

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 151.033905
$ws.Range("H15").Value = 453.101715
$ws.Range("I15").Value = 0.1782978134070997
$ws.Range("J15").Value = 0.1782978134070997
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.6655859999999999
$ws.Range("N15").Value = 1.996758
$ws.Range("O15").Value = 0.1138523176430159
$ws.Range("P15").Value = 0.1138523176430159
$ws.Range("Q15").Value = 100.52605269333
$ws.Range("R15").Value = 904.7344742399698
$ws.Range("S15").Value = 0.02029961928708029
$ws.Range("T15").Value = 0.02029961928708029

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 151.033905
$ws.Range("H16").Value = 453.101715
$ws.Range("I16").Value = 0.1782978134070997
$ws.Range("J16").Value = 0.1782978134070997
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.01401
$ws.Range("N16").Value = 3.04203
$ws.Range("O16").Value = 0.1734522490154458
$ws.Range("P16").Value = 0.1734522490154458
$ws.Range("Q16").Value = 153.14989000905
$ws.Range("R16").Value = 1378.34901008145
$ws.Range("S16").Value = 0.03092615672999776
$ws.Range("T16").Value = 0.03092615672999775

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 130.8644813333333
$ws.Range("H17").Value = 392.593444
$ws.Range("I17").Value = 0.1544875031496243
$ws.Range("J17").Value = 0.1544875031496243
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.166450999999999
$ws.Range("N17").Value = 12.499353
$ws.Range("O17").Value = 0.7126954333415383
$ws.Range("P17").Value = 0.7126954333415383
$ws.Range("Q17").Value = 545.2404491157479
$ws.Range("R17").Value = 4907.164042041732
$ws.Range("S17").Value = 0.1101025380030737
$ws.Range("T17").Value = 0.1101025380030737

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 130.8644813333333
$ws.Range("H18").Value = 392.593444
$ws.Range("I18").Value = 0.1544875031496243
$ws.Range("J18").Value = 0.1544875031496243
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 0.6655859999999999
$ws.Range("N18").Value = 1.996758
$ws.Range("O18").Value = 0.1138523176430159
$ws.Range("P18").Value = 0.1138523176430159
$ws.Range("Q18").Value = 87.10156667272797
$ws.Range("R18").Value = 783.9141000545519
$ws.Range("S18").Value = 0.01758876028046743
$ws.Range("T18").Value = 0.01758876028046743

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 130.8644813333333
$ws.Range("H19").Value = 392.593444
$ws.Range("I19").Value = 0.1544875031496243
$ws.Range("J19").Value = 0.1544875031496243
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 1.01401
$ws.Range("N19").Value = 3.04203
$ws.Range("O19").Value = 0.1734522490154458
$ws.Range("P19").Value = 0.1734522490154458
$ws.Range("Q19").Value = 132.6978927168133
$ws.Range("R19").Value = 1194.28103445132
$ws.Range("S19").Value = 0.0267962048660831
$ws.Range("T19").Value = 0.02679620486608309

